$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data as of the latest GitHub Actions run.
# Row 40 and 41 swap coins (ARBITRUM <-> WEMIXToken) in addition to value updates.

$ws.Range('D2').Value = '26.966.73'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.556.98'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').Value = "'206.88"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').Value = "'0.489"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').Value = "'22.07"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = "'0.0595"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '1.779.86'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = '1.505.30'
$ws.Range('E13').Value = '  -3.03%  '
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = "'0.521"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').Value = "'61.94"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '26.971.40'
$ws.Range('D18').Value = '0.0₃0707'
$ws.Range('E18').Value = '  +2.71%  '
$ws.Range('D19').Value = "'217.25"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('D20').Value = "'7.32"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('E22').Value = '  +1.39%  '
$ws.Range('D23').Value = "'9.23"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.02%  '
$ws.Range('E24').Value = '  -3.73%  '
$ws.Range('D25').Value = "'152.87"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = "'6.64"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').Value = "'15.01"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('D30').Value = "'0.0469"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').Value = "'3.22"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = '1.422.12'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').Value = "'3.11"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.68%  '
$ws.Range('D35').Value = "'1.09"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.56%  '
$ws.Range('D36').Value = "'1.60"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.80%  '
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').Value = "'0.528"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').Value = "'1.03"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.49%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = "'0.807"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.30%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  +1.96%  '
$ws.Range('E44').Value = '  +2.29%  '
$ws.Range('D45').Value = "'64.80"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('D47').Value = '1.693.20'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = "'87.44"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('D49').Value = "'0.0520"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').Value = '0.0₇0995'
$ws.Range('E50').Value = '  +0.81%  '
$ws.Range('D51').Value = "'0.0957"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.09%  '
